$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must keep General/text
# semantics identical to the source data (plain text price strings), so
# force Text number format before writing those specific cells.

$ws.Range("D2").Value = "37.524.18"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.078.38"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.32"
$ws.Range("E7").Value = "  +5.58%  "
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.88"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("D13").Value = "2.383.47"
$ws.Range("E13").Value = "  +3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.56"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.13"
$ws.Range("E15").Value = "  +3.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.781"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.21"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "2.063.51"
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("D19").Value = "37.685.72"
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  +17.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.96"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.56"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.47"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.50"
$ws.Range("E28").Value = "  +9.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.03"
$ws.Range("E29").Value = "  +3.82%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.29"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.53"
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0626"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  +7.79%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.59"
$ws.Range("E41").Value = "  +19.48%  "
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0960"
$ws.Range("E43").Value = "  +3.90%  "
$ws.Range("D44").Value = "1.475.40"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +6.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "95.93"
$ws.Range("E46").Value = "  +5.65%  "
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.83"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("E50").Value = "  +5.49%  "
$ws.Range("E51").Value = "  +1.56%  "
